$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New cells C7/D7 (Step3 row for Example1 bank2/bank3 block) ----
# C7/D7 reuse the "quote-prefix" text styles of the rows above (styles 14/15),
# so copy formatting first, then fill the value using a leading apostrophe so
# Excel keeps treating the cell as quote-prefixed text (not a formula).
$ws.Range("C6").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("C7").Value = "'Step3"

$ws.Range("D6").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("D7").Value = "'=bank2!=bank3"

# ---- New block starting at row 53 - Example2 header/table ----
$ws.Range("B53:C53").Merge()
$ws.Range("C3").Copy()
$ws.Range("B53:C53").PasteSpecial(-4122)
$ws.Range("B53").Value = "Spreadsheet SpreadsheetResult  Example2 (Integer[] intArr1, int[] intArr2, Double[] dArr, Float[] floatArr, int i1, Integer i2, float f)"

$ws.Range("C4").Copy()
$ws.Range("B54").PasteSpecial(-4122)
$ws.Range("B54").Value = "Step"

$ws.Range("D4").Copy()
$ws.Range("C54").PasteSpecial(-4122)
$ws.Range("C54").Value = "Formula"

$ws.Range("C5").Copy()
$ws.Range("B55").PasteSpecial(-4122)
$ws.Range("B55").Value = "'Step1"

$ws.Range("D5").Copy()
$ws.Range("C55").PasteSpecial(-4122)
$ws.Range("C55").Value = "'=intArr1>i1"

$ws.Range("C5").Copy()
$ws.Range("B56").PasteSpecial(-4122)
$ws.Range("B56").Value = "'Step2"

$ws.Range("D5").Copy()
$ws.Range("C56").PasteSpecial(-4122)
$ws.Range("C56").Value = "'=intArr2<i2"

$ws.Range("C5").Copy()
$ws.Range("B57").PasteSpecial(-4122)
$ws.Range("B57").Value = "'Step3"

$ws.Range("D5").Copy()
$ws.Range("C57").PasteSpecial(-4122)
$ws.Range("C57").Value = "'=dArr>=i1"

$ws.Range("C5").Copy()
$ws.Range("B58").PasteSpecial(-4122)
$ws.Range("B58").Value = "'Step4"

$ws.Range("D5").Copy()
$ws.Range("C58").PasteSpecial(-4122)
$ws.Range("C58").Value = "'=dArr<=floatArr"

$ws.Range("C5").Copy()
$ws.Range("B59").PasteSpecial(-4122)
$ws.Range("B59").Value = "'Step5"

$ws.Range("D5").Copy()
$ws.Range("C59").PasteSpecial(-4122)
$ws.Range("C59").Value = "'=floatArr<=f"
